# t13.2: fix the "Trimestre" date strings (01/10/YYYY -> 31/10/YYYY) for the
# existing rows and append the newly published 31/10/2024 quarter rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Correct the existing Trimestre column (C2:C97): 01/10/YYYY -> 31/10/YYYY
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $old = [string]$cell.Value2
    $new = $old -replace '^01/10/', '31/10/'
    $cell.Value = $new
}

# --- 2) Append the new 31/10/2024 quarter rows (98:105)
$newRows = @(
    @("Sergipe", "Agricultura, pecuária, produção florestal, pesca e aquicultura", "31/10/2024", 116, 11.4),
    @("Sergipe", "Indústria geral", "31/10/2024", 82, 8),
    @("Sergipe", "Construção", "31/10/2024", 83, 8.1),
    @("Sergipe", "Comércio, reparação de veículos automotores e motocicletas", "31/10/2024", 200, 19.7),
    @("Sergipe", "Transporte, armazenagem e correio", "31/10/2024", 50, 4.9),
    @("Sergipe", "Alojamento e alimentação", "31/10/2024", 58, 5.7),
    @("Sergipe", "Informação, comunicação e atividades financeiras, imobiliárias, profissionais e administrativas", "31/10/2024", 102, 10),
    @("Sergipe", "Administração pública, defesa, seguridade social, educação, saúde humana e serviços sociais", "31/10/2024", 209, 20.5)
)

$rowIndex = 98
foreach ($row in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
